# Insert a new data row before the existing row 174 (shifts rows 174-266
# down to 175-267) and populate it with the new "Ají" price observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("174").Insert()

$ws.Range("A174").Value = 4
$ws.Range("B174").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C174").Value = "Los Lagos"
$ws.Range("D174").Value = 44719
$ws.Range("E174").Value = 10
$ws.Range("F174").Value = 100112021
$ws.Range("G174").Value = "Ají"
$ws.Range("H174").Value = "Inferno"
$ws.Range("I174").Value = "Primera"
$ws.Range("J174").Value = 80
$ws.Range("K174").Value = 30000
$ws.Range("L174").Value = 30000
$ws.Range("M174").Value = 30000
$ws.Range("N174").Value = '$/caja 12 kilos'
$ws.Range("O174").Value = "Región de Arica y Parinacota"
$ws.Range("P174").Value = 2500
$ws.Range("Q174").Value = 12
$ws.Range("R174").Value = "Hortaliza"
